# Remove the trailing "Ver no Jupiter ..." and "(c) 2020 ..." paragraphs,
# along with the blank paragraph that separates them from the
# "LOB1049: ..." requirement line, while leaving the blank paragraph and
# page-break paragraph that follow them untouched.

$d = $word.ActiveDocument

$jupiterText   = "Ver no Jupiter Salvar em pdf Salvar em docx"
$copyrightText = [char]0x00A9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

# Walk paragraphs from the end so deleting one doesn't disturb the index
# of paragraphs we still need to inspect/delete.
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]0x0007, [char]0x000D, [char]0x000A)

    if ($text -eq $copyrightText) {
        $para.Range.Delete()
        continue
    }

    if ($text -eq $jupiterText) {
        # The blank paragraph immediately preceding this one (the one that
        # sits right after the "LOB1049..." requirement paragraph) is also
        # removed as part of this edit.
        $prev = $para.Previous()
        $para.Range.Delete()
        if ($prev -ne $null) {
            $prevText = $prev.Range.Text.TrimEnd([char]0x0007, [char]0x000D, [char]0x000A)
            if ($prevText -eq "") {
                $prev.Range.Delete()
            }
        }
        continue
    }
}
